$p = $ppt.ActivePresentation
$s = $p.Slides.Item(30)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Find the paragraph ending in "Can the store learn to do profitable business with Cluster 2-4?"
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i)
    $text = $para.Text.TrimEnd("`r")
    if ($text -eq "Can the store learn to do profitable business with Cluster 2-4?") {
        $para.InsertAfter("`rGet more granular data`rGet more data on profitability rather than just revenue") | Out-Null
        break
    }
}
